$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title: "Play Ritual Respins Free: Voodoo-themed Slot Review"
#    -> "Play Ritual Respins for Free"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Ritual Respins Free: Voodoo-themed Slot Review",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Ritual Respins for Free", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the whole "Meta description: ..." paragraph entirely.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3. "What we like" bullet list changes
# ------------------------------------------------------------------

# 3a. Insert a new bullet "Fun Voodoo theme" right before
#     "Smooth gameplay on any device".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Smooth gameplay on any device*") {
        $p.Range.InsertParagraphBefore()
        $newp = $d.Paragraphs($i)
        $newp.Range.Text = "Fun Voodoo theme"
        break
    }
}

# 3b. "Low volatility with frequent wins" -> "243 ways to win"
$d.Content.Find.Execute(
    "Low volatility with frequent wins",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "243 ways to win", 2) | Out-Null

# 3c. Remove the "Free spins with sticky bonus symbols" bullet entirely.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Free spins with sticky bonus symbols*") {
        $p.Range.Delete()
        break
    }
}

# 3d. "Fun theme and light atmosphere" -> "Low volatility"
$d.Content.Find.Execute(
    "Fun theme and light atmosphere",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Low volatility", 2) | Out-Null

# ------------------------------------------------------------------
# 4. "What we don't like" bullet list changes
# ------------------------------------------------------------------

# 4a. "Second Scatter needed to win during free spins" -> "Not high payouts"
$d.Content.Find.Execute(
    "Second Scatter needed to win during free spins",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Not high payouts", 2) | Out-Null

# 4b. Insert a new (non-bulleted, bold) paragraph right after
#     the "Not high payouts" bullet: "Play Ritual Respins for Free"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Not high payouts*") {
        $p.Range.InsertParagraphAfter()
        $newp = $d.Paragraphs($i + 1)
        $newp.Style = "Normal"
        $r = $newp.Range
        $startPos = $r.Start
        $r.Text = "Play Ritual Respins for Free"
        $boldRange = $d.Range($startPos, $startPos + 29)
        $boldRange.Font.Bold = $true
        break
    }
}

# ------------------------------------------------------------------
# 5. Final (italic) paragraph: replace the image-prompt text with the
#    new short review blurb.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for " + [char]34 + "Ritual Respins" + [char]34 + " that features a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior holding a cauldron of magical potion and surrounded by symbols from the game, such as playing cards, voodoo dolls, and potions. There should be a sense of fun and adventure in the image, with the Maya warrior looking like he is ready for an exciting game of slots. The background could be a witch's lair, with candles and magic books, to set the atmosphere of the game. Overall, the image should make viewers want to try their luck with the game and see what magical bonuses and winning opportunities await them.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Ritual Respins, a fun Voodoo-themed slot game with 243 ways to win.", 2) | Out-Null
